$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SkillLevelConfig.attributeType used to be stored as a free-form lowercase
# string ("attack"). Switch it over to the new AttributeType enum member
# names (upper snake-case), per-row:
#   row 5 (level 1) -> ATTACK
#   row 6 (level 2) -> ATTACK
#   row 7 (level 3) -> CRITICAL_RATE
$ws.Range("H5").Value = "ATTACK"
$ws.Range("H6").Value = "ATTACK"
$ws.Range("H7").Value = "CRITICAL_RATE"

# Leave the cursor where the author last left it when saving.
$ws.Range("G13").Select()
